$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.072.53"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.302.73"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'300.22"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'98.03"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'17.69"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "2.661.18"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "2.282.78"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "'0.788"
$ws.Range("D18").Value = "42.970.27"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'12.68"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'68.57"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "'237.96"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "'25.04"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "'164.36"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("E30").Value = "  -12.92%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'33.05"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.11"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "'18.09"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").Value = "'2.41"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.0697"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "2.021.08"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").Value = "'10.36"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "'17.48"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "'2.84"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").Value = "'54.26"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "2.528.00"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  -0.87%  "
